$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("H2").Value = 82
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 5

# Update the active selection/cell to H2
$ws.Range("H2").Select()
